$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows (10-15) entirely, keeping 1-9
$ws.Range("A10:G15").EntireRow.Delete()

# Row 2
$ws.Range("B2").Value = "B_High"
$ws.Range("C2").Value = "AT221"
$ws.Range("D2").Value = "Centralized"
$ws.Range("E2").Value = "GWh/km**2"
$ws.Range("F2").Value = 4.3899999999999997
$ws.Range("G2").Value = "Techno-Friendly"

# Row 3
$ws.Range("B3").Value = "A_Low"
$ws.Range("C3").Value = "AT221"
$ws.Range("D3").Value = "Centralized"
$ws.Range("E3").Value = "GWh/km**2"
$ws.Range("F3").Value = 1.75
$ws.Range("G3").Value = "Directed Transition"

# Row 4
$ws.Range("B4").Value = "B_High"
$ws.Range("C4").Value = "AT312"
$ws.Range("D4").Value = "Centralized"
$ws.Range("E4").Value = "GWh/km**2"
$ws.Range("F4").Value = 4.9800000000000004
$ws.Range("G4").Value = "Societal Commitment"

# Row 5
$ws.Range("B5").Value = "A_Low"
$ws.Range("C5").Value = "AT312"
$ws.Range("D5").Value = "Centralized"
$ws.Range("E5").Value = "GWh/km**2"
$ws.Range("F5").Value = 2.7
$ws.Range("G5").Value = "Gradual Development"

# Row 6
$ws.Range("B6").Value = "B_High"
$ws.Range("C6").Value = "AT342"
$ws.Range("D6").Value = "Centralized"
$ws.Range("E6").Value = "GWh/km**2"
$ws.Range("F6").Value = 2.13
$ws.Range("G6").Value = "Techno-Friendly"

# Row 7
$ws.Range("B7").Value = "A_Low"
$ws.Range("C7").Value = "AT342"
$ws.Range("D7").Value = "Centralized"
$ws.Range("E7").Value = "GWh/km**2"
$ws.Range("F7").Value = 1.65
$ws.Range("G7").Value = "Directed Transition"

# Row 8
$ws.Range("B8").Value = "B_High"
$ws.Range("C8").Value = "AT130"
$ws.Range("D8").Value = "Centralized"
$ws.Range("E8").Value = "GWh/km**2"
$ws.Range("F8").Value = 16.62
$ws.Range("G8").Value = "Gradual Development"

# Row 9
$ws.Range("B9").Value = "A_Low"
$ws.Range("C9").Value = "AT130"
$ws.Range("D9").Value = "Centralized"
$ws.Range("E9").Value = "GWh/km**2"
$ws.Range("F9").Value = 10.17
$ws.Range("G9").Value = "Directed Transition"

# Column width adjustments (best-fit sized to new, shorter content)
$ws.Columns.Item(1).ColumnWidth = 12.21875
$ws.Columns.Item(2).ColumnWidth = 8

# View: zoom + selection
$excel.ActiveWindow.Zoom = 175
$ws.Range("D15").Select()
